# v1.0.10 - Add analyses for all three MCP1, MCP1_pg_ml_2015, and MCP1_pg_ug_2015.
# Add comparison between MCP1, MCP1_pg_ml_2015, and MCP1_pg_ug_2015.
# Add (and fixed) ordinal regression. Double checked which measurement to use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (columns A-L), in final post-edit order:
# Dataset | Outcome | CpG | Beta | s.e.m. | HR | low95CI | up95CI | Z-value | P-value | SampleSize | N_events
$rowsData = @(
    @("AEDB.CEA", "epmajor.90days", "MCP1_pg_ug_2015_rank", 0.917425873382322, 0.326801161762416, 2.50283946506284, 1.31902865241112, 4.74910486320781, 2.80729073432514, 0.00499601286701933, 1186, 47),
    @("AEDB.CEA", "epmajor.90days", "MCP1_pg_ml_2015_rank", 0.548431292434352, 0.32837567051412, 1.73053618226251, 0.909204696239631, 3.29381875226302, 1.67013375739957, 0.0948929025340179, 1187, 47),
    @("AEDB.CEA", "epmajor.90days", "MCP1_rank", -0.45028416451512, 0.464233230895538, 0.637446986068871, 0.256614249324403, 1.58346101636238, -0.969952460418419, 0.332070189030053, 549, 21),
    @("AEDB.CEA", "epstroke.90days", "MCP1_pg_ug_2015_rank", 0.865632589223619, 0.366771002231371, 2.3765089645354, 1.15807764152256, 4.87687064840645, 2.36014456965589, 0.0182678144077115, 1186, 36),
    @("AEDB.CEA", "epstroke.90days", "MCP1_pg_ml_2015_rank", 0.701347068344952, 0.372875805194614, 2.01646719788885, 0.970941062918279, 4.18783396382522, 1.8809133190578, 0.0599837111064671, 1187, 36),
    @("AEDB.CEA", "epstroke.90days", "MCP1_rank", -0.12295996907087, 0.51929575950932, 0.884299061165344, 0.319569769628737, 2.44699250021799, -0.236782155099926, 0.812825796101932, 549, 16),
    @("AEDB.CEA", "epcoronary.90days", "MCP1_pg_ug_2015_rank", 0.564743098833264, 0.548834786095088, 1.75899583658885, 0.59991153163458, 5.15753772011435, 1.02898561305008, 0.303486435501202, 1186, 16),
    @("AEDB.CEA", "epcoronary.90days", "MCP1_pg_ml_2015_rank", 0.475207927418731, 0.570195625833423, 1.60834858248985, 0.526041327881079, 4.91745614972271, 0.833412088568983, 0.404612359312045, 1187, 16),
    @("AEDB.CEA", "epcoronary.90days", "MCP1_rank", -0.493062384133497, 0.789443614057766, 0.610753166570324, 0.129980366391938, 2.86981365594008, -0.624569475708518, 0.532253658644671, 549, 7),
    @("AEDB.CEA", "epcvdeath.90days", "MCP1_pg_ug_2015_rank", 0.266979011415606, 0.795670031647141, 1.30601303467819, 0.274574061145918, 6.21205819526736, 0.335539860490818, 0.737217881345752, 1186, 8),
    @("AEDB.CEA", "epcvdeath.90days", "MCP1_pg_ml_2015_rank", -1.38472537610728, 0.901589847253064, 0.250392554128394, 0.0427732882848567, 1.46578469126344, -1.5358706404317, 0.124570110132459, 1187, 8),
    @("AEDB.CEA", "epcvdeath.90days", "MCP1_rank", -0.349307143289704, 1.24019098590833, 0.705176506772318, 0.0620334529543105, 8.01622160336403, -0.281655928206789, 0.778207350992465, 549, 3)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
